# Apply the pedalboard-display BoM update:
#  - JST PH connector rows get a real datasheet link (was "~") and a Digikey
#    supplier link (was empty), on both the BoM sheet and the DNF sheet.
#  - Those rows grow taller (30pt) to accommodate the wrapped text.

$wb = $excel.ActiveWorkbook

$datasheetUrl = "https://www.jst-mfg.com/product/pdf/eng/ePH.pdf"
$supplier3pin = "https://www.digikey.ch/de/products/detail/jst-sales-america-inc/B3B-PH-SM4-TB/926832"
$supplier4pin = "https://www.digikey.ch/de/products/detail/jst-sales-america-inc./B4B-PH-SM4-TB/926833"

# Fill colors (OLE BGR-encoded values) matching the existing "status" fills
# already used in the workbook's style table, so Excel's style de-dup reuses
# the same cellXf rather than fabricating a new one.
$colorStyle7  = 11790079   # FFE6B3 - orange  (cellXfs index 7)
$colorStyle8  = 16775654   # E6F9FF - lt blue (cellXfs index 8)
$colorStyle11 = 12447999   # FFF0BD - yellow  (cellXfs index 11)
$colorStyle12 = 16777200   # F0FFFF - lt cyan (cellXfs index 12)

# ---- BoM sheet ----
$bom = $wb.Worksheets.Item("BoM")

$bom.Range("I11").Value = $datasheetUrl
$bom.Range("I11").Interior.Color = $colorStyle7
$bom.Range("J11").Value = $supplier3pin
$bom.Range("J11").Interior.Color = $colorStyle8
$bom.Rows.Item(11).RowHeight = 30

$bom.Range("I12").Value = $datasheetUrl
$bom.Range("I12").Interior.Color = $colorStyle11
$bom.Range("J12").Value = $supplier4pin
$bom.Range("J12").Interior.Color = $colorStyle12
$bom.Rows.Item(12).RowHeight = 30

# ---- DNF sheet ----
$dnf = $wb.Worksheets.Item("DNF")

$dnf.Range("I9").Value = $datasheetUrl
$dnf.Range("I9").Interior.Color = $colorStyle7
$dnf.Range("J9").Value = $supplier3pin
$dnf.Range("J9").Interior.Color = $colorStyle8
$dnf.Rows.Item(9).RowHeight = 30
